# chore(results): update lottery results 2025-09-19T17:41:23Z
#
# Appends the latest "Pick 3" draw as a new row at the bottom of the
# Results sheet, matching how every prior row in this sheet was written
# (plain text values for every column, including the numeric-looking
# Phase/Result codes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet view is left-to-right (explicit, matching the rest of the export).
$ws.DisplayRightToLeft = $false

# Locate the first empty row below the existing data.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$rowRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 5))

# Date ("2025-09-19") and Phase ("250919") both look numeric/date-like to
# Excel's auto-detection, so force the whole row to Text first to keep
# every column a plain string, exactly like the existing rows.
$rowRange.NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "2025-09-19"
$ws.Cells.Item($newRow, 2).Value = "Pick 3"
$ws.Cells.Item($newRow, 3).Value = "250919"
$ws.Cells.Item($newRow, 4).Value = "5-9-2"
$ws.Cells.Item($newRow, 5).Value = "2025-09-19T21:41:22.972+04:00"

# Drop the temporary "@" text formatting so the new row keeps the same
# (unformatted / default-style) look as the rows above it - only the
# stored text values are meant to differ, not the cell formatting.
$rowRange.ClearFormats()
